$wb = $excel.ActiveWorkbook
$wsValid = $wb.Worksheets.Item("Valid Species")
$wsMain = $wb.Worksheets.Item("Main")

# The "Valid Species" sheet is protected; unprotect it so we can edit the list.
$wsValid.Unprotect()

# --- Add the new species "Pogoniulus atroflavus" ---
# The list is kept in alphabetical order, so insert it right before the
# existing "Pogoniulus bilineatus" row (found dynamically in case row
# numbers shift).
$anchorAdd = $wsValid.Columns(1).Find("Pogoniulus bilineatus")
$addRow = $anchorAdd.Row
$wsValid.Rows($addRow + ":" + $addRow).Insert()
$wsValid.Cells.Item($addRow, 1).Value = "Pogoniulus atroflavus"

# --- Remove "Turdus philomelos" and all of its subspecies rows ---
$firstRemove = $wsValid.Columns(1).Find("Turdus philomelos")
$lastRemove = $wsValid.Columns(1).Find("Turdus philomelos philomelos")
$startRow = $firstRemove.Row
$endRow = $lastRemove.Row
$wsValid.Rows($startRow + ":" + $endRow).Delete()

# Restore protection on the "Valid Species" sheet.
$wsValid.Protect()

# The data validation list on "Main"!A2:A1048576 points at the "Valid
# Species" range; update its upper bound to match the sheet's new row count.
$lastRow = $wsValid.UsedRange.Rows.Count
$dv = $wsMain.Range("A2").Validation
$dv.Formula1 = "='Valid Species'!A1:A" + $lastRow
